$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2: close trade 1 ---
$ws.Range("F2").Value = "CLOSED"
$ws.Range("G2").Value = "TP"
$ws.Range("H2").Value = 112
$ws.Range("I2").Value = 1268
$ws.Range("M2").Value = "TP"

# --- Add row 4: trade 3, already closed ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "USDJPY"
$ws.Range("C4").Value = "sell"
$ws.Range("D4").Value = 67
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = "CLOSED"
$ws.Range("G4").Value = "SL"
$ws.Range("H4").Value = -67
$ws.Range("I4").Value = 1201
$ws.Range("J4").Value = 45727.40790072917
$ws.Range("J4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = "SL"
